# Update NATMI ligand-receptor edge statistics for Sema6b-Plxna2
# following the revised pipeline (Dr Hou advice): ligand- and
# receptor-expressing cell counts go from 1 to 3, which changes the
# downstream expression/specificity/edge-weight figures accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.893586666666666
$ws.Range("H2").Value = 26.68076
$ws.Range("I2").Value = 0.7066310027692895
$ws.Range("J2").Value = 0.7066310027692896
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 29.420614
$ws.Range("N2").Value = 88.261842
$ws.Range("O2").Value = 0.5865186809777162
$ws.Range("P2").Value = 0.5865186809777162
$ws.Range("Q2").Value = 261.6547803955467
$ws.Range("R2").Value = 2354.89302355992
$ws.Range("S2").Value = 0.4144522836822046
$ws.Range("T2").Value = 0.4144522836822047

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.893586666666666
$ws.Range("H3").Value = 26.68076
$ws.Range("I3").Value = 0.7066310027692895
$ws.Range("J3").Value = 0.7066310027692896
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.080312
$ws.Range("N3").Value = 12.240936
$ws.Range("O3").Value = 0.08134361887272465
$ws.Range("P3").Value = 0.08134361887272466
$ws.Range("Q3").Value = 36.28860839904
$ws.Range("R3").Value = 326.59747559136
$ws.Range("S3").Value = 0.05747992297291632
$ws.Range("T3").Value = 0.05747992297291634

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.893586666666666
$ws.Range("H4").Value = 26.68076
$ws.Range("I4").Value = 0.7066310027692895
$ws.Range("J4").Value = 0.7066310027692896
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.660501
$ws.Range("N4").Value = 49.981503
$ws.Range("O4").Value = 0.3321377001495591
$ws.Range("P4").Value = 0.3321377001495591
$ws.Range("Q4").Value = 148.1716095535866
$ws.Range("R4").Value = 1333.54448598228
$ws.Range("S4").Value = 0.2346987961141686
$ws.Range("T4").Value = 0.2346987961141686

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.082649333333334
$ws.Range("H5").Value = 6.247948000000001
$ws.Range("I5").Value = 0.1654748125799407
$ws.Range("J5").Value = 0.1654748125799407
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.420614
$ws.Range("N5").Value = 88.261842
$ws.Range("O5").Value = 0.5865186809777162
$ws.Range("P5").Value = 0.5865186809777162
$ws.Range("Q5").Value = 61.27282213335735
$ws.Range("R5").Value = 551.4553992002161
$ws.Range("S5").Value = 0.0970540688094216
$ws.Range("T5").Value = 0.0970540688094216

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.082649333333334
$ws.Range("H6").Value = 6.247948000000001
$ws.Range("I6").Value = 0.1654748125799407
$ws.Range("J6").Value = 0.1654748125799407
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.080312
$ws.Range("N6").Value = 12.240936
$ws.Range("O6").Value = 0.08134361887272465
$ws.Range("P6").Value = 0.08134361887272466
$ws.Range("Q6").Value = 8.497859066592003
$ws.Range("R6").Value = 76.48073159932802
$ws.Range("S6").Value = 0.01346032008753824
$ws.Range("T6").Value = 0.01346032008753824

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.082649333333334
$ws.Range("H7").Value = 6.247948000000001
$ws.Range("I7").Value = 0.1654748125799407
$ws.Range("J7").Value = 0.1654748125799407
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.660501
$ws.Range("N7").Value = 49.981503
$ws.Range("O7").Value = 0.3321377001495591
$ws.Range("P7").Value = 0.3321377001495591
$ws.Range("Q7").Value = 34.69798130064934
$ws.Range("R7").Value = 312.2818317058441
$ws.Range("S7").Value = 0.05496042368298083
$ws.Range("T7").Value = 0.05496042368298083

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.609663333333333
$ws.Range("H8").Value = 4.82899
$ws.Range("I8").Value = 0.1278941846507698
$ws.Range("J8").Value = 0.1278941846507698
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.420614
$ws.Range("N8").Value = 88.261842
$ws.Range("O8").Value = 0.5865186809777162
$ws.Range("P8").Value = 0.5865186809777162
$ws.Range("Q8").Value = 47.35728359995333
$ws.Range("R8").Value = 426.21555239958
$ws.Range("S8").Value = 0.07501232848608996
$ws.Range("T8").Value = 0.07501232848608996

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.609663333333333
$ws.Range("H9").Value = 4.82899
$ws.Range("I9").Value = 0.1278941846507698
$ws.Range("J9").Value = 0.1278941846507698
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.080312
$ws.Range("N9").Value = 12.240936
$ws.Range("O9").Value = 0.08134361887272465
$ws.Range("P9").Value = 0.08134361887272466
$ws.Range("Q9").Value = 6.56792861496
$ws.Range("R9").Value = 59.11135753464001
$ws.Range("S9").Value = 0.01040337581227009
$ws.Range("T9").Value = 0.01040337581227009

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.609663333333333
$ws.Range("H10").Value = 4.82899
$ws.Range("I10").Value = 0.1278941846507698
$ws.Range("J10").Value = 0.1278941846507698
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.660501
$ws.Range("N10").Value = 49.981503
$ws.Range("O10").Value = 0.3321377001495591
$ws.Range("P10").Value = 0.3321377001495591
$ws.Range("Q10").Value = 26.81779757466333
$ws.Range("R10").Value = 241.36017817197
$ws.Range("S10").Value = 0.04247848035240972
$ws.Range("T10").Value = 0.04247848035240972
